# 1) Save Password work implemented. 2) Form Authentication redirect work done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8: "Role based authentication" now assigned to Shambhoo ---
$ws.Range("C8").Value = "Shambhoo"

# --- Row 16: "Save Password" -> reassigned to Vipin, marked Done ---
$ws.Range("C16").Value = "Vipin"
$ws.Range("D16").Value = "Done"

# --- Row 19: "Form authentication... redirect" work -> marked Done ---
$ws.Range("D19").Value = "Done"

# --- Row 34: "Dependency Injection" -> back to Pending, remark cleared ---
$ws.Range("D34").Value = "Pending"
$ws.Range("G34").Value = ""

# --- Row 35: now holds the "invalid user message" task ---
$ws.Range("B35").Value = "In case of any error, login page still show invalid user message"

# --- Row 36 (new task row): Login background image change ---
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "Login background image change"
$ws.Range("C36").Value = "Vipin"
$ws.Range("D36").Value = "Pending"
$ws.Range("E36").Value = "Account"

# --- Row 37 (new task row): Sign Out ---
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "Sign Out"
$ws.Range("C37").Value = "Vipin"
$ws.Range("D37").Value = "Done"
$ws.Range("E37").Value = "Account"

# --- Row 38 (new task row): Put logged in user detail in session ---
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "Put logged in user detail in session"
$ws.Range("C38").Value = "Shambhoo"
$ws.Range("D38").Value = "Pending"
$ws.Range("E38").Value = "Account"

# --- Row 23: remark now needs a decision from Shambhoo/Vipin ---
$ws.Range("C23").Value = "Shambhoo/Vipin"

# --- Row 33: add remark "Need to decide with Suraj" ---
$ws.Range("G33").Value = "Need to decide with Suraj"
$ws.Range("G33").WrapText = $true

# --- Extend the autofilter range so it covers the new rows (A1:G35 -> A1:G38) ---
# Toggle off then back on over the new range; this refreshes the stored
# <autoFilter ref="..."> without introducing any live filter criteria
# (so no rows get hidden as a side effect).
$ws.Range("A1:G38").AutoFilter()
$ws.Range("A1:G38").AutoFilter()

# --- Keep the hidden _FilterDatabase defined name in sync with the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$38"
    }
}
